# Remove the trailing blank slide (slide 8) from the deck.
# This slide has an empty shape tree and is the last slide in the
# presentation (p:sldId id="264" r:id="rId9" in presentation.xml).
$p = $ppt.ActivePresentation
$p.Slides.Item($p.Slides.Count).Delete()
